$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows at the top of the data (right after the header row),
# pushing all existing data rows down by 6.
$ws.Rows("2:7").Insert()

# The newly inserted rows are blank/unstyled; clone the date / percent
# number formats from row 8 (the first row that still carries the
# original formatting, since it held the data that used to be row 2).
$ws.Range("A8").Copy()
$ws.Range("A2:A7").PasteSpecial(-4122)
$ws.Range("D8").Copy()
$ws.Range("D2:D7").PasteSpecial(-4122)

# Populate the 6 brand-new days of data (most recent first).
$ws.Range("A2").Value = 45560
$ws.Range("B2").Value = -19.25
$ws.Range("C2").Value = -321.89
$ws.Range("D2").Value = -0.13669999999999999

$ws.Range("A3").Value = 45559
$ws.Range("B3").Value = 20.9
$ws.Range("C3").Value = -302.64999999999998
$ws.Range("D3").Value = -0.1285

$ws.Range("A4").Value = 45558
$ws.Range("B4").Value = 31.1
$ws.Range("C4").Value = -323.55
$ws.Range("D4").Value = -0.13739999999999999

$ws.Range("A5").Value = 45557
$ws.Range("B5").Value = -23.78
$ws.Range("C5").Value = -354.64
$ws.Range("D5").Value = -0.15060000000000001

$ws.Range("A6").Value = 45556
$ws.Range("B6").Value = 17.09
$ws.Range("C6").Value = -330.86
$ws.Range("D6").Value = -0.14050000000000001

$ws.Range("A7").Value = 45555
$ws.Range("B7").Value = 25.21
$ws.Range("C7").Value = -347.95
$ws.Range("D7").Value = -0.14779999999999999

# Row 8 already holds what used to be row 2's data (and its formatting);
# overwrite with the recalculated values for that date.
$ws.Range("B8").Value = 32.979999999999997
$ws.Range("C8").Value = -373.17
$ws.Range("D8").Value = -0.1585

# Column width / autofit tweaks recorded in the saved file
# (ColumnWidth values chosen to land on the closest width this engine stores).
$ws.Columns("C").ColumnWidth = 13.6666666666667
$ws.Columns("D").ColumnWidth = 7

# Selection moved to C2 before saving.
$ws.Range("C2").Select() | Out-Null
